# Addition of PDMS contaminant series to database
# Silicon added as basic elemental component and polydimethylsiloxane
# oligomer contaminant series added to LOBDbase.

$wb = $excel.ActiveWorkbook
$wsRt = $wb.Worksheets.Item(1)     # "LOBSTAHS_rt.windows"
$wsNotes = $wb.Worksheets.Item(2)  # "Notes"

# --- 1. Log the change in the "Notes" history table (row 26) -------------
# Copy formatting from the previous history row (25) down into row 26,
# then fill in the new entry: date, description, initials.
$wsNotes.Range("A25:C25").Copy()
[void]$wsNotes.Range("A26:C26").PasteSpecial(-4122)
$wsNotes.Range("A26").Value = 42779
$wsNotes.Range("B26").Value = "Added PDMS"
$wsNotes.Range("C26").Value = "JEH"

# --- 2. Append the new PDMS6-PDMS27 contaminant rows to the RT table -----
$pdmsNames = @("PDMS6","PDMS7","PDMS8","PDMS9","PDMS10","PDMS11","PDMS12", `
  "PDMS13","PDMS14","PDMS15","PDMS16","PDMS17","PDMS18","PDMS19","PDMS20", `
  "PDMS21","PDMS22","PDMS23","PDMS24","PDMS25","PDMS26","PDMS27")

$row = 73
foreach ($name in $pdmsNames) {
    $wsRt.Range("A$row").Value = $name
    $wsRt.Range("A$row").HorizontalAlignment = -4152   # xlRight
    $wsRt.Range("B$row").Value = 30
    $wsRt.Range("C$row").Value = 5
    $row = $row + 1
}

# --- 3. Restore cursor / active-sheet state seen in the saved workbook ---
[void]$wsNotes.Range("B27").Select()
$wsRt.Activate()
[void]$wsRt.Range("F82").Select()
